$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy H1's formatting (bold, borders, centered) onto I1/J1, then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-22: I = 1 (constant), J = same as H
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
